$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update: new spheroid measurements were inserted as the first
# --- data point of the multi_poly_* block (cols I:L), pushing the
# --- existing rows 24-30 down into rows 25-31 (the old row-31 values,
# --- which had nothing left to push into, drop off the bottom).
$ws.Range("I24").Value = 96.7
$ws.Range("J24").Value = 65
$ws.Range("K24").Value = 60.9
$ws.Range("L24").Value = 768

$ws.Range("I25").Value = 46.4
$ws.Range("J25").Value = 31.8
$ws.Range("K25").Value = 31
$ws.Range("L25").Value = 64.400000000000006

$ws.Range("I26").Value = 83.4
$ws.Range("J26").Value = 67.8
$ws.Range("K26").Value = 59.6
$ws.Range("L26").Value = 442.5

$ws.Range("I27").Value = 58.6
$ws.Range("J27").Value = 54.6
$ws.Range("K27").Value = 51.8
$ws.Range("L27").Value = 210

$ws.Range("I28").Value = 97.3
$ws.Range("J28").Value = 65.599999999999994
$ws.Range("K28").Value = 54
$ws.Range("L28").Value = 594.1

$ws.Range("I29").Value = 69.099999999999994
$ws.Range("J29").Value = 48.9
$ws.Range("K29").Value = 45.9
$ws.Range("L29").Value = 223.7

$ws.Range("I30").Value = 91.1
$ws.Range("J30").Value = 63.3
$ws.Range("K30").Value = 45.8
$ws.Range("L30").Value = 371.5

$ws.Range("I31").Value = 47.2
$ws.Range("J31").Value = 36.6
$ws.Range("K31").Value = 34.700000000000003
$ws.Range("L31").Value = 113.8

# --- Row 91 had four formatted-but-empty cells (E91:H91) left over from a
# --- wider table; drop them entirely so the row matches the later rows.
$ws.Range("E91:H91").Clear()

# --- Move the cursor / selection to where the author left it.
$ws.Range("Q13").Select()
